# Encryption/decryption test results workbook - add the "Results" summary
# table (row/column headers + AVERAGE/MAX/MIN stats over Table1) next to the
# existing data, and slide the two charts over to make room for it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary table -----------------------------------------------------
# Row labels first (G22/G23), then the column headers (H21/I21/J21) - this
# mirrors the order the strings were actually typed in, which is also the
# order they land in the shared-strings table.
$ws.Range("G22").Value2 = "Encryption Time"
$ws.Range("G23").Value2 = "Decryption Time"

$ws.Range("H21").Value2 = "Average (ms)"
$ws.Range("I21").Value2 = "Max (ms)"
$ws.Range("J21").Value2 = "Min (ms)"

# Stats pulled from the Table1 structured reference columns.
$ws.Range("H22").Formula = "=AVERAGE(Table1[Encryption time (ms)])"
$ws.Range("I22").Formula = "=MAX(Table1[Encryption time (ms)])"
$ws.Range("J22").Formula = "=MIN(Table1[Encryption time (ms)])"

$ws.Range("H23").Formula = "=AVERAGE(Table1[Successful decryption time (ms)])"
$ws.Range("I23").Formula = "=MAX(Table1[Successful decryption time (ms)])"
$ws.Range("J23").Formula = "=MIN(Table1[Successful decryption time (ms)])"

# --- Reposition the two charts to make room for the new table ----------
# (done before resizing columns G/H/J below, since the anchor cell/offset
# the engine recomputes on save depends on the column widths at this time)
$chartObjs = $ws.ChartObjects()

$co1 = $chartObjs.Item(1)
$co1.Left = 1044.9559759473425
$co1.Top = 10.87503937007874
$co1.Width = 706.5
$co1.Height = 204.37496062992125

$co2 = $chartObjs.Item(2)
$co2.Left = 1222.8935153174214
$co2.Top = 227.62496062992125
$co2.Width = 443.5
$co2.Height = 216.0

# --- Column widths for the new table (best-fit to the header/label text) -
$ws.Columns.Item(7).ColumnWidth = 23.834
$ws.Columns.Item(8).ColumnWidth = 11.834
$ws.Columns.Item(10).ColumnWidth = 8.001

# --- Selection matches where the user left off --------------------------
$ws.Range("G21:J23").Select()
